$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.494.89'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '1.619.39'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'211.09"
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = "'22.82"
$ws.Range('E8').Value = '  -1.08%  '
$ws.Range('D9').Value = "'0.263"
$ws.Range('E9').Value = '  +2.09%  '
$ws.Range('D10').Value = "'0.0613"
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').Value = '1.848.28'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').Value = '1.613.70'
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').Value = "'4.03"
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').Value = "'0.549"
$ws.Range('E15').Value = '  -2.52%  '
$ws.Range('D16').Value = "'65.09"
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '27.490.28'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = "'230.32"
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('E20').Value = '  -1.55%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'4.28"
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D23').Value = "'10.18"
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('E24').Value = '  +6.08%  '
$ws.Range('D25').Value = "'149.90"
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = "'6.83"
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.111"
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').Value = "'15.57"
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').Value = "'0.0483"
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  -1.16%  '
$ws.Range('D33').Value = '1.448.57'
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('E34').Value = '  -3.43%  '
$ws.Range('E35').Value = '  -3.57%  '
$ws.Range('D36').Value = "'2.33"
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').Value = "'0.934"
$ws.Range('E37').Value = '  +4.12%  '
$ws.Range('D38').Value = "'0.561"
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('D40').Value = "'0.863"
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = "'67.95"
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('B43').Value = 'mCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D43').Value = "'2.49"
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'0.990"
$ws.Range('E44').Value = '  -3.78%  '
$ws.Range('D45').Value = "'5.43"
$ws.Range('E45').Value = '  -4.34%  '
$ws.Range('D47').Value = '1.759.23'
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').Value = "'86.53"
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('E50').Value = '  +17.86%  '
$ws.Range('E51').Value = '  +1.57%  '
